$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.954.79"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "'1.818.40"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("E7").Value = "  +0.43%  "
$ws.Range("D8").Value = "'0.3663"
$ws.Range("E8").Value = "  -1.10%  "
$ws.Range("D9").Value = "'0.07356"
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("D10").Value = "'0.8734"
$ws.Range("E10").Value = "  -0.46%  "
$ws.Range("D11").Value = "'20.26"
$ws.Range("E11").Value = "  -1.07%  "
$ws.Range("D12").Value = "'1.819.59"
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("D13").Value = "'5.415"
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("D14").Value = "'0.07113"
$ws.Range("E14").Value = "  +0.99%  "
$ws.Range("D15").Value = "'6.516"
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("D16").Value = "'91.49"
$ws.Range("E16").Value = "  -0.27%  "
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").Value = "'0.000008706"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").Value = "'14.67"
$ws.Range("E20").Value = "  -0.60%  "
$ws.Range("D21").Value = "'26.974.76"
$ws.Range("D22").Value = "'5.294"
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("D23").Value = "'10.60"
$ws.Range("E23").Value = "  +0.35%  "
$ws.Range("D24").Value = "'2.051.95"
$ws.Range("E24").Value = "  +0.78%  "
$ws.Range("D25").Value = "'1.890"
$ws.Range("E25").Value = "  -0.57%  "
$ws.Range("D26").Value = "'150.64"
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("D27").Value = "'18.42"
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("D28").Value = "'2.140"
$ws.Range("E28").Value = "  -0.46%  "
$ws.Range("D29").Value = "'5.251"
$ws.Range("E29").Value = "  -1.45%  "
$ws.Range("D30").Value = "'116.60"
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("D31").Value = "'0.08895"
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("D32").Value = "'0.7596"
$ws.Range("E32").Value = "  +0.73%  "
$ws.Range("D33").Value = "'1.163"
$ws.Range("E33").Value = "  +0.77%  "
$ws.Range("D34").Value = "'4.505"
$ws.Range("E34").Value = "  +0.90%  "
$ws.Range("E35").Value = "  -0.36%  "
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("D37").Value = "'1.094"
$ws.Range("E37").Value = "  -0.52%  "
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("D39").Value = "'0.01947"
$ws.Range("E39").Value = "  -0.97%  "
$ws.Range("D40").Value = "'2.970"
$ws.Range("E40").Value = "  +1.53%  "
$ws.Range("D41").Value = "'7.183"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").Value = "'0.5293"
$ws.Range("E42").Value = "  -0.62%  "
$ws.Range("D43").Value = "'2.343"
$ws.Range("E43").Value = "  -3.90%  "
$ws.Range("D44").Value = "'0.1654"
$ws.Range("E44").Value = "  -0.68%  "
$ws.Range("D45").Value = "'8.447"
$ws.Range("E45").Value = "  -0.48%  "
$ws.Range("D46").Value = "'0.4869"
$ws.Range("E46").Value = "  -2.12%  "
$ws.Range("D47").Value = "'10.48"
$ws.Range("E47").Value = "  +1.03%  "
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("D49").Value = "'1.668"
$ws.Range("E49").Value = "  -0.28%  "
$ws.Range("D50").Value = "'103.35"
$ws.Range("E50").Value = "  -0.29%  "
$ws.Range("E51").Value = "  +0.03%  "
